# Update betting odds values for rows 3, 4, 6, and 15 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 2.88
$ws.Range("I3").Value = 4.33
$ws.Range("R3").Value = 2.2
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 5.5
$ws.Range("U3").Value = 8.5
$ws.Range("W3").Value = 17
$ws.Range("AA3").Value = 6
$ws.Range("AB3").Value = 19

# Row 4
$ws.Range("G4").Value = 3.6
$ws.Range("I4").Value = 2.2
$ws.Range("L4").Value = 1.5
$ws.Range("M4").Value = 2.5
$ws.Range("O4").Value = 1.47
$ws.Range("P4").Value = 1.57
$ws.Range("Q4").Value = 2.25
$ws.Range("S4").Value = 1.63
$ws.Range("V4").Value = 13
$ws.Range("AB4").Value = 19
$ws.Range("AF4").Value = 10
$ws.Range("AG4").Value = 21

# Row 6
$ws.Range("O6").Value = 1.63
$ws.Range("R6").Value = 1.87
$ws.Range("S6").Value = 1.77

# Row 15
$ws.Range("G15").Value = 1.8
$ws.Range("I15").Value = 3.3
$ws.Range("K15").Value = 29
$ws.Range("L15").Value = 1.1
$ws.Range("M15").Value = 7
$ws.Range("T15").Value = 15
$ws.Range("U15").Value = 13
$ws.Range("W15").Value = 19
$ws.Range("X15").Value = 13
$ws.Range("Z15").Value = 29
$ws.Range("AA15").Value = 10
$ws.Range("AE15").Value = 23
$ws.Range("AF15").Value = 12
$ws.Range("AH15").Value = 23
$ws.Range("AI15").Value = 21
